# The template used a real Word field (fldChar begin/instrText/fldChar end)
# to hold the M2Doc "m:null.yesNo()" query. The parser was updated to use
# TokenIteratorFieldRewriterSplit, which expects the query to be written as
# plain literal text delimited by "{" and "}" instead of a Word field.
#
# Find the paragraph that still contains the old field code and rewrite its
# content as plain-text runs: "{" "m" ":" "null" "." "yesNo" "(" ")" "}"
# keeping the existing bookmark and the orange run-color on "null".

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if ($r.Fields.Count -gt 0) {
        $newXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
            '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
            '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData>' +
            '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
            '<w:body>' +
            '<w:p>' +
            '<w:r><w:t>{</w:t></w:r>' +
            '<w:r><w:t>m</w:t></w:r>' +
            '<w:r><w:t>:</w:t></w:r>' +
            '<w:r><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>null</w:t></w:r>' +
            '<w:r><w:t>.</w:t></w:r>' +
            '<w:r><w:t>yesNo</w:t></w:r>' +
            '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
            '<w:bookmarkEnd w:id="0"/>' +
            '<w:r><w:t>()</w:t></w:r>' +
            '<w:r><w:t xml:space="preserve">}</w:t></w:r>' +
            '</w:p>' +
            '</w:body>' +
            '</w:document>' +
            '</pkg:xmlData>' +
            '</pkg:part>' +
            '</pkg:package>'
        $r.InsertXML($newXml)
        break
    }
}
